# Fruta / hortaliza, semanal
# Insert a new weekly record as row 13, pushing the existing rows 13-42
# down to 14-43 (the sheet keeps the same repeating set of records, just
# shifted by one position to make room for the newest entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 13..42 down to 14..43, leaving a blank row 13 to fill in.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the latest weekly price record.
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44708
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100114007
$ws.Range("G13").Value = "Jengibre"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 440
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13500
$ws.Range("N13").Value = "$/caja 13 kilos"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 1038
$ws.Range("Q13").Value = 13
$ws.Range("R13").Value = "Hortaliza"
